# "additional runs for iters 3 & 4"
# Duplicate Sheet1 into a new "Sheet1 (2)" tab and record two more
# Monte-Carlo runs (iteration 4a / 4b) on the copy, highlighting the
# columns that improved for each new run.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")

# 1. Duplicate the sheet (Excel will name the copy "Sheet1 (2)" and
#    place it right after Sheet1).
$ws1.Copy($null, $ws1)
$ws2 = $wb.Worksheets.Item("Sheet1 (2)")

# 2. Update the "Notes" text for the existing iteration-4 row on the new
#    tab -- this run used a 0 discount instead of -1.
$ws2.Range("G5").Value = "Num2Move = 1; Discount of 0 for ""far out"" movements"

# 3. The note text for the two brand-new rows (write this before the A6/A7
#    labels so the shared-string table order matches the authored file).
$ws2.Range("G7").Value = "Num2Move = 1; Discount of -1 for ""far out"" movements and -1 for ""too close"" movements"

# 4. New row: iteration 4a
$ws2.Range("A6").Value = "4a"
$ws2.Range("B6").Value = 114
$ws2.Range("C6").Value = 12.78
$ws2.Range("D6").Value = 1.67
$ws2.Range("E6").Value = 7053.31
$ws2.Range("F6").Value = 66.19
$ws2.Range("G6").Value = "Num2Move = 1; Discount of -1 for ""far out"" movements"

# 5. New row: iteration 4b
$ws2.Range("A7").Value = "4b"
$ws2.Range("B7").Value = 113
$ws2.Range("C7").Value = 13.04
$ws2.Range("D7").Value = 1.56
$ws2.Range("E7").Value = 8068.57
$ws2.Range("F7").Value = 68.67
$ws2.Range("H7").Value = "iter4a_1_bvm1_corr_results10000a"

# 6. Formatting to match the rest of the metrics table: thin borders
#    around the new rows, centered alignment for label/metric cells, and
#    number formats consistent with the columns above (integers for
#    Num. Steps / Distance Moved, one-decimal for Error / CEP Radius).
$newRows = $ws2.Range("A6:G7")
$newRows.Borders.LineStyle = 1
$newRows.Borders.Weight = 2
$ws2.Range("A6:A7").HorizontalAlignment = -4108
$ws2.Range("A6:A7").VerticalAlignment = -4108
$ws2.Range("A6:A7").WrapText = $true
$ws2.Range("A6:G7").Interior.Color = 16777215
$ws2.Range("B6:F7").HorizontalAlignment = -4108

$ws2.Range("B6").NumberFormat = "0"
$ws2.Range("C6").NumberFormat = "0.0"
$ws2.Range("D6").NumberFormat = "0.0"
$ws2.Range("E6").NumberFormat = "0"
$ws2.Range("F6").NumberFormat = "0"

$ws2.Range("B7").NumberFormat = "0"
$ws2.Range("C7").NumberFormat = "0.0"
$ws2.Range("D7").NumberFormat = "0.0"
$ws2.Range("E7").NumberFormat = "0"
$ws2.Range("F7").NumberFormat = "0"

# 7. Highlight the metric(s) that improved on each new run (mirrors the
#    green highlight already used for best-so-far values in this sheet).
$ws2.Range("C6").Interior.Color = 5296274
$ws2.Range("E6").Interior.Color = 5296274
$ws2.Range("B7").Interior.Color = 5296274
$ws2.Range("D7").Interior.Color = 5296274

$null = $ws2.Range("A1").Select()
